# Generate Report for Handoff
# A new file (8d6db238-c65f-4a1d-bb81-137e821782eb) has been handed off.
# On every sheet (Overview, zh-cn, de-de) a new row is inserted as row 2
# holding the new file's data, and the previous row 2 (f6befd66-...) is
# pushed down to row 3, unchanged.

$wb = $excel.ActiveWorkbook

$HL_UNDERLINE = 2         # xlUnderlineStyleSingle
$HL_COLOR = 15570276      # RGB(100,149,237) == FF6495ED, BGR-encoded for COM

function Restore-HyperlinkLook($rng) {
    $rng.Font.Name = "Calibri"
    $rng.Font.Size = 11
    $rng.Font.Underline = $HL_UNDERLINE
    $rng.Font.Color = $HL_COLOR
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Duplicate row 2 (old data) down into row 3, keeping its formatting intact.
$ws.Range("A2:D2").Copy()
$ws.Range("A3").Insert(-4121)

# Clear every hyperlink (all currently still refer to row 2's old file) and
# rebuild them for both rows in their final form.
$ws.Hyperlinks.Delete()

# Row 2 -> brand-new handoff entry
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9c7d1a24f5b6e8c3d2a1f0e9b8c7d6a5f4e3d2c1/e2e/8d6db238-c65f-4a1d-bb81-137e821782eb.md", "", "", "8d6db238-c65f-4a1d-bb81-137e821782eb.md")
Restore-HyperlinkLook $ws.Range("A2")
$ws.Range("B2").Value = "Ready for handoff"
$ws.Range("C2").Value = "Ready for handoff"
$ws.Range("D2").Value = "2016-06-17 10:06:55"

# Row 3 -> previous entry, unchanged
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8baeb8214d3d5d0c059a20fdc49c7166d060299c/e2e/f6befd66-f219-45b9-9d97-52234355658b.md", "", "", "f6befd66-f219-45b9-9d97-52234355658b.md")
Restore-HyperlinkLook $ws.Range("A3")
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-06-17 10:06:36"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2:K2").Copy()
$ws.Range("A3").Insert(-4121)

$ws.Hyperlinks.Delete()

# Row 2 -> brand-new handoff entry
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9c7d1a24f5b6e8c3d2a1f0e9b8c7d6a5f4e3d2c1/e2e/8d6db238-c65f-4a1d-bb81-137e821782eb.md", "", "", "8d6db238-c65f-4a1d-bb81-137e821782eb.md")
Restore-HyperlinkLook $ws.Range("A2")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/9c7d1a24f5b6e8c3d2a1f0e9b8c7d6a5f4e3d2c1/e2e/8d6db238-c65f-4a1d-bb81-137e821782eb.md", "", "", ".md")
Restore-HyperlinkLook $ws.Range("B2")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1a2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8d6db238-c65f-4a1d-bb81-137e821782eb.095c687f8f34ad07ef6111827f0c73cf68b48cd8.zh-cn.xlf", "", "", "8d6db238-c65f-4a1d-bb81-137e821782eb.095c687f8f34ad07ef6111827f0c73cf68b48cd8.zh-cn.xlf")
Restore-HyperlinkLook $ws.Range("D2")
$ws.Range("E2").Value = "2016-03-17 10:06:51"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

# Row 3 -> previous entry, unchanged
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8baeb8214d3d5d0c059a20fdc49c7166d060299c/e2e/f6befd66-f219-45b9-9d97-52234355658b.md", "", "", "f6befd66-f219-45b9-9d97-52234355658b.md")
Restore-HyperlinkLook $ws.Range("A3")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/8baeb8214d3d5d0c059a20fdc49c7166d060299c/e2e/f6befd66-f219-45b9-9d97-52234355658b.md", "", "", ".md")
Restore-HyperlinkLook $ws.Range("B3")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ccafc8bd75387903e10189b3f0c3d03086e84943/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f6befd66-f219-45b9-9d97-52234355658b.d3566084f2f22f9b736bdd9db7eb607a8413b37a.zh-cn.xlf", "", "", "f6befd66-f219-45b9-9d97-52234355658b.d3566084f2f22f9b736bdd9db7eb607a8413b37a.zh-cn.xlf")
Restore-HyperlinkLook $ws.Range("D3")
$ws.Range("E3").Value = "2016-03-17 10:06:32"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2:K2").Copy()
$ws.Range("A3").Insert(-4121)

$ws.Hyperlinks.Delete()

# Row 2 -> brand-new handoff entry
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/9c7d1a24f5b6e8c3d2a1f0e9b8c7d6a5f4e3d2c1/e2e/8d6db238-c65f-4a1d-bb81-137e821782eb.md", "", "", "8d6db238-c65f-4a1d-bb81-137e821782eb.md")
Restore-HyperlinkLook $ws.Range("A2")
$ws.Hyperlinks.Add($ws.Range("B2"), "https://github.com/OpenLocalizationTest/oltest/blob/9c7d1a24f5b6e8c3d2a1f0e9b8c7d6a5f4e3d2c1/e2e/8d6db238-c65f-4a1d-bb81-137e821782eb.md", "", "", ".md")
Restore-HyperlinkLook $ws.Range("B2")
$ws.Range("C2").Value = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2b3c4d5e6f7a8b9c0d1e2f3a4b5c6d7e8f9a0b1c/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8d6db238-c65f-4a1d-bb81-137e821782eb.095c687f8f34ad07ef6111827f0c73cf68b48cd8.de-de.xlf", "", "", "8d6db238-c65f-4a1d-bb81-137e821782eb.095c687f8f34ad07ef6111827f0c73cf68b48cd8.de-de.xlf")
Restore-HyperlinkLook $ws.Range("D2")
$ws.Range("E2").Value = "2016-03-17 10:06:55"
$ws.Range("H2").Value = "0001-01-01 00:00:00"
$ws.Range("I2").Value = "Include"

# Row 3 -> previous entry, unchanged
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/8baeb8214d3d5d0c059a20fdc49c7166d060299c/e2e/f6befd66-f219-45b9-9d97-52234355658b.md", "", "", "f6befd66-f219-45b9-9d97-52234355658b.md")
Restore-HyperlinkLook $ws.Range("A3")
$ws.Hyperlinks.Add($ws.Range("B3"), "https://github.com/OpenLocalizationTest/oltest/blob/8baeb8214d3d5d0c059a20fdc49c7166d060299c/e2e/f6befd66-f219-45b9-9d97-52234355658b.md", "", "", ".md")
Restore-HyperlinkLook $ws.Range("B3")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Hyperlinks.Add($ws.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c14541241ed22e7eb5503b9539980edef1e87ee8/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f6befd66-f219-45b9-9d97-52234355658b.d3566084f2f22f9b736bdd9db7eb607a8413b37a.de-de.xlf", "", "", "f6befd66-f219-45b9-9d97-52234355658b.d3566084f2f22f9b736bdd9db7eb607a8413b37a.de-de.xlf")
Restore-HyperlinkLook $ws.Range("D3")
$ws.Range("E3").Value = "2016-03-17 10:06:36"
$ws.Range("H3").Value = "0001-01-01 00:00:00"
$ws.Range("I3").Value = "Include"
